# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732; G = 6.201049113329182 }
    3 = @{ B = 0.6753301551942219; C = 0.3127903958511391;  D = 0.1575252929769615;  E = 0.496779210170732; G = 1.642425054193055 }
    4 = @{ B = 0.6753301551942219; C = 0.3127903958511391;  D = 0.1575252929769615;  E = 0.496779210170732; G = 1.642425054193055 }
    5 = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732; G = 4.429675500412797 }
    6 = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732; G = 4.429675500412797 }
    7 = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732; G = 4.429675500412797 }
    8 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 26.21740644021617;   E = 0.496779210170732; G = 31.61296591696135 }
    9 = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732; G = 4.429675500412797 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
